$d = $word.ActiveDocument

# The last paragraph in the body currently holds only the _GoBack bookmark.
# We replace that paragraph's range with: two new paragraphs (a spacer, and
# the bold/underlined question), an empty paragraph, the plain-text answer
# paragraph, and finally a reconstructed paragraph carrying the original
# bookmark so the document structure matches the target edit.
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $last.Range

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="60" w:after="60" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="3D464D"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:lang w:eastAsia="pt-BR"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="60" w:after="60" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:color w:val="3D464D"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:u w:val="single"/><w:lang w:eastAsia="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:color w:val="3D464D"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:u w:val="single"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Você comprou 5 peras por R$ 1,50 cada e precisa calcular o valor total e terá que fazer o cálculo para demais produtos também. Qual a vantagem de utilizar o Excel ao invés de sua calculadora manual, se você precisar calcular o valor total de vários produtos adquiridos, conforme a quantidade comprada?</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:color w:val="3D464D"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Através do Excel, você poderá criar uma fórmula única, referenciando a célula que contém o valor unitário e multiplicando pela célula que contém a quantidade. Criada a fórmula, você pode propagar o cálculo para as demais linhas, copiando a fórmula da célula principal. Já na calculadora manual, você teria que calcular produto a produto.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

[void]$r.InsertXML($xml)
